$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.02270466666666667
$ws.Range("H2").Value = 0.06811400000000001
$ws.Range("I2").Value = 0.002206225855740089
$ws.Range("J2").Value = 0.002206225855740089
$ws.Range("M2").Value = 4.861952666666666
$ws.Range("N2").Value = 14.585858
$ws.Range("O2").Value = 0.3995648519435639
$ws.Range("P2").Value = 0.3995648519435638
$ws.Range("Q2").Value = 0.1103890146457778
$ws.Range("R2").Value = 0.9935011318119999
$ws.Range("S2").Value = 0.0008815303074028511
$ws.Range("T2").Value = 0.0008815303074028509
$ws.Range("G3").Value = 0.02270466666666667
$ws.Range("H3").Value = 0.06811400000000001
$ws.Range("I3").Value = 0.002206225855740089
$ws.Range("J3").Value = 0.002206225855740089
$ws.Range("O3").Value = 0.04932556406896855
$ws.Range("P3").Value = 0.04932556406896854
$ws.Range("Q3").Value = 0.01362732579688889
$ws.Range("R3").Value = 0.122645932172
$ws.Range("S3").Value = 0.0001088233347979227
$ws.Range("T3").Value = 0.0001088233347979227
$ws.Range("G4").Value = 0.02270466666666667
$ws.Range("H4").Value = 0.06811400000000001
$ws.Range("I4").Value = 0.002206225855740089
$ws.Range("J4").Value = 0.002206225855740089
$ws.Range("M4").Value = 4.206754333333333
$ws.Range("N4").Value = 12.620263
$ws.Range("O4").Value = 0.3457193616641432
$ws.Range("P4").Value = 0.3457193616641432
$ws.Range("Q4").Value = 0.09551295488688888
$ws.Range("R4").Value = 0.859616593982
$ws.Range("S4").Value = 0.0007627349945333917
$ws.Range("T4").Value = 0.0007627349945333917
$ws.Range("G5").Value = 0.02270466666666667
$ws.Range("H5").Value = 0.06811400000000001
$ws.Range("I5").Value = 0.002206225855740089
$ws.Range("J5").Value = 0.002206225855740089
$ws.Range("M5").Value = 2.499212666666667
$ws.Range("N5").Value = 7.497638
$ws.Range("O5").Value = 0.2053902223233243
$ws.Range("P5").Value = 0.2053902223233243
$ws.Range("Q5").Value = 0.05674379052577778
$ws.Range("R5").Value = 0.5106941147320001
$ws.Range("S5").Value = 0.0004531372190059233
$ws.Range("T5").Value = 0.0004531372190059233
$ws.Range("I6").Value = 0.002281111990432972
$ws.Range("J6").Value = 0.002281111990432972
$ws.Range("M6").Value = 4.861952666666666
$ws.Range("N6").Value = 14.585858
$ws.Range("O6").Value = 0.3995648519435639
$ws.Range("P6").Value = 0.3995648519435638
$ws.Range("Q6").Value = 0.1141359595008889
$ws.Range("R6").Value = 1.027223635508
$ws.Range("S6").Value = 0.0009114521747240389
$ws.Range("T6").Value = 0.0009114521747240388
$ws.Range("I7").Value = 0.002281111990432972
$ws.Range("J7").Value = 0.002281111990432972
$ws.Range("O7").Value = 0.04932556406896855
$ws.Range("P7").Value = 0.04932556406896854
$ws.Range("S7").Value = 0.0001125171356325939
$ws.Range("T7").Value = 0.0001125171356325939
$ws.Range("I8").Value = 0.002281111990432972
$ws.Range("J8").Value = 0.002281111990432972
$ws.Range("M8").Value = 4.206754333333333
$ws.Range("N8").Value = 12.620263
$ws.Range("O8").Value = 0.3457193616641432
$ws.Range("P8").Value = 0.3457193616641432
$ws.Range("Q8").Value = 0.09875496022644444
$ws.Range("R8").Value = 0.8887946420379998
$ws.Range("S8").Value = 0.0007886245812169105
$ws.Range("T8").Value = 0.0007886245812169105
$ws.Range("I9").Value = 0.002281111990432972
$ws.Range("J9").Value = 0.002281111990432972
$ws.Range("M9").Value = 2.499212666666667
$ws.Range("N9").Value = 7.497638
$ws.Range("O9").Value = 0.2053902223233243
$ws.Range("P9").Value = 0.2053902223233243
$ws.Range("Q9").Value = 0.0586698504208889
$ws.Range("R9").Value = 0.528028653788
$ws.Range("S9").Value = 0.0004685180988594291
$ws.Range("T9").Value = 0.0004685180988594291
$ws.Range("G10").Value = 10.24499966666667
$ws.Range("H10").Value = 30.734999
$ws.Range("I10").Value = 0.9955126621538269
$ws.Range("J10").Value = 0.9955126621538269
$ws.Range("M10").Value = 4.861952666666666
$ws.Range("N10").Value = 14.585858
$ws.Range("O10").Value = 0.3995648519435639
$ws.Range("P10").Value = 0.3995648519435638
$ws.Range("Q10").Value = 49.81070344934911
$ws.Range("R10").Value = 448.296331044142
$ws.Range("S10").Value = 0.397771869461437
$ws.Range("T10").Value = 0.3977718694614369
$ws.Range("G11").Value = 10.24499966666667
$ws.Range("H11").Value = 30.734999
$ws.Range("I11").Value = 0.9955126621538269
$ws.Range("J11").Value = 0.9955126621538269
$ws.Range("O11").Value = 0.04932556406896855
$ws.Range("P11").Value = 0.04932556406896854
$ws.Range("Q11").Value = 6.149041969933556
$ws.Range("R11").Value = 55.34137772940201
$ws.Range("S11").Value = 0.04910422359853803
$ws.Range("T11").Value = 0.04910422359853803
$ws.Range("G12").Value = 10.24499966666667
$ws.Range("H12").Value = 30.734999
$ws.Range("I12").Value = 0.9955126621538269
$ws.Range("J12").Value = 0.9955126621538269
$ws.Range("M12").Value = 4.206754333333333
$ws.Range("N12").Value = 12.620263
$ws.Range("O12").Value = 0.3457193616641432
$ws.Range("P12").Value = 0.3457193616641432
$ws.Range("Q12").Value = 43.09819674274855
$ws.Range("R12").Value = 387.8837706847369
$ws.Range("S12").Value = 0.3441680020883929
$ws.Range("T12").Value = 0.3441680020883929
$ws.Range("G13").Value = 10.24499966666667
$ws.Range("H13").Value = 30.734999
$ws.Range("I13").Value = 0.9955126621538269
$ws.Range("J13").Value = 0.9955126621538269
$ws.Range("M13").Value = 2.499212666666667
$ws.Range("N13").Value = 7.497638
$ws.Range("O13").Value = 0.2053902223233243
$ws.Range("P13").Value = 0.2053902223233243
$ws.Range("Q13").Value = 25.60443293692911
$ws.Range("R13").Value = 230.439896432362
$ws.Range("S13").Value = 0.204468567005459
$ws.Range("T13").Value = 0.204468567005459
